$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Anshal) ---
# Total no. of Posts: 17 -> 19
$ws.Range("C2").Value = "'19"
# No. of times Reported: 1 -> 3
$ws.Range("D2").Value = "'3"

# --- Update existing row 5 (f20190363) ---
# Total no. of Posts: 0 -> 1
$ws.Range("C5").Value = "'1"
# D5 and E5 remain unchanged (0 / False)

# --- Append new row 6 (sambhav) ---
$ws.Range("A6").Value = "sambhav"
$ws.Range("B6").Value = "f20190192@pilani.bits-pilani.ac.in"
$ws.Range("C6").Value = "'0"
$ws.Range("D6").Value = "'0"
$ws.Range("E6").Value = "'False"
